$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitplanung")

# Row 25 (Anforderung #01 / 307) - add hours in KW 51 (Freitag column Y)
$ws.Range("Y25").Value = 1

# Row 26 (Anforderung #02 / 308) - add hours in KW 51 (Y) and KW 52 (AB)
$ws.Range("Y26").Value = 2
$ws.Range("AB26").Value = 1

# Row 27 (Anforderung #03 / 309) - add hours in KW 1 (AR, AS)
$ws.Range("AR27").Value = 2
$ws.Range("AS27").Value = 2

# Row 32 (Testfälle erstellen / 401) - add hours in KW 51 (Y), KW 52 (AB, AC), KW 1 (AR, AS)
$ws.Range("Y32").Value = 1.5
$ws.Range("AB32").Value = 3
$ws.Range("AC32").Value = 2
$ws.Range("AR32").Value = 2
$ws.Range("AS32").Value = 2

# Row 33 (Bugs fixen / 402) - add hours in KW 51 (Y), KW 52 (AB, AC), KW 1 (AR, AS)
$ws.Range("Y33").Value = 1
$ws.Range("AB33").Value = 2
$ws.Range("AC33").Value = 2.5
$ws.Range("AR33").Value = 3
$ws.Range("AS33").Value = 2

# Row 34 (Testfälle durchführen + dokumentieren / 403) - add hours in KW 52 (AB, AC), KW 1 (AR, AS)
$ws.Range("AB34").Value = 1.5
$ws.Range("AC34").Value = 3
$ws.Range("AR34").Value = 0.5
$ws.Range("AS34").Value = 1.5

# Row 40 (601) - add hours in KW 1 (AT)
$ws.Range("AT40").Value = 2

# Row 41 (602) - add hours in KW 1 (AT)
$ws.Range("AT41").Value = 2

# Row 42 (603) - new task "PowerPoint" with hours in KW 1 (AT)
$ws.Range("B42").Value = "PowerPoint"
$ws.Range("AT42").Value = 3.5

# Scroll the sheet view so row 6 is at the top and select T46 (matches the
# author's view state when the workbook was saved)
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("T46").Select() | Out-Null
